# Applies the KPI recomputation described in the commit message:
# "Change the computations of the KPIs"
#
# Three sheets are touched:
#   - Productdata:     column C (StartingInventories) and column E (SetupCosts)
#                       for rows 2..18
#   - Capacity:        column B for rows 2..18
#   - ProcessingTime:  the diagonal non-zero entries for rows 2..17

$wb = $excel.ActiveWorkbook

# --- Productdata sheet: StartingInventories (C) / SetupCosts (E) ---
$wsProd = $wb.Worksheets.Item("Productdata")

$startingInventories = @{
    2  = 5
    3  = 5
    4  = 5
    5  = 5
    6  = 0
    7  = 0
    8  = 0
    9  = 0
    10 = 0
    11 = 0
    12 = 0
    13 = 0
    14 = 0
    15 = 0
    16 = 0
    17 = 0
    18 = 0
}

$setupCosts = @{
    2  = 1.2375
    3  = 1.84635
    4  = 1.8729
    5  = 1.91133
    6  = 1.08
    7  = 3.1725
    8  = 1.77885
    9  = 3.5208
    10 = 0.495
    11 = 1.08
    12 = 0.72
    13 = 1.485
    14 = 2.025
    15 = 1.08
    16 = 0.67635
    17 = 1.3518
    18 = 3.1725
}

foreach ($row in 2..18) {
    $wsProd.Range("C$row").Value = $startingInventories[$row]
    $wsProd.Range("E$row").Value = $setupCosts[$row]
}

# --- Capacity sheet: column B ---
$wsCap = $wb.Worksheets.Item("Capacity")

$capacity = @{
    2  = 30
    3  = 20
    4  = 30
    5  = 50
    6  = 30
    7  = 30
    8  = 50
    9  = 100
    10 = 40
    11 = 200
    12 = 200
    13 = 150
    14 = 120
    15 = 20
    16 = 50
    17 = 60
    18 = 120
}

foreach ($row in 2..18) {
    $wsCap.Range("B$row").Value = $capacity[$row]
}

# --- ProcessingTime sheet: diagonal entries ---
$wsProc = $wb.Worksheets.Item("ProcessingTime")

$processingTime = @{
    "B2"  = 3
    "C3"  = 2
    "E5"  = 5
    "F6"  = 3
    "G7"  = 1
    "H8"  = 5
    "I9"  = 5
    "L12" = 5
    "M13" = 5
    "N14" = 4
    "O15" = 2
    "P16" = 5
    "Q17" = 3
}

foreach ($cellRef in $processingTime.Keys) {
    $wsProc.Range($cellRef).Value = $processingTime[$cellRef]
}
